# Fixed courseDb type representation & add export xls
#
# 1. The "Type" column (D) used the abbreviation "L" for every class row;
#    spell it out as "lecture".
# 2. The "Sem" column (A) on the data rows was rendering with the theme's
#    automatic text color instead of plain black - force it to solid black
#    (RGB 0,0,0 / FF000000) like the rest of the numeric columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "L" class-type abbreviation to "lecture" -------------
$ws.Range("D2").Value = "lecture"
$ws.Range("D3").Value = "lecture"
$ws.Range("D4").Value = "lecture"

# --- 2. Make the Sem column's font solid black on the data rows ---------
$ws.Range("A2").Font.Color = 0
$ws.Range("A3").Font.Color = 0
$ws.Range("A4").Font.Color = 0
